$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "final" with "true " (trailing space) in column E
$ws.Range("E5").Value = "true "
$ws.Range("E7").Value = "true "
$ws.Range("E9").Value = "true "
$ws.Range("E11").Value = "true "

# Apply strikethrough font style to F9:I9
$ws.Range("F9:I9").Font.Strikethrough = $true

# Update selection to C8
$ws.Range("C8").Select()
